$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "65.469.83"
Set-TextValue "E2" "  -3.38%  "

# Row 3
Set-TextValue "D3" "3.490.01"
Set-TextValue "E3" "  -0.71%  "

# Row 4
Set-TextValue "D4" "0.999"
Set-TextValue "E4" "  -0.20%  "

# Row 5
Set-TextValue "D5" "553.54"
Set-TextValue "E5" "  -0.64%  "

# Row 6
Set-TextValue "D6" "179.54"
Set-TextValue "E6" "  -6.16%  "

# Row 7
Set-TextValue "D7" "0.641"
Set-TextValue "E7" "  +4.68%  "

# Row 8
Set-TextValue "E8" "  -0.10%  "

# Row 9
Set-TextValue "D9" "0.633"
Set-TextValue "E9" "  -0.95%  "

# Row 10
Set-TextValue "D10" "0.156"
Set-TextValue "E10" "  +3.16%  "

# Row 11
Set-TextValue "D11" "53.96"
Set-TextValue "E11" "  -5.60%  "

# Row 12
Set-TextValue "D12" "0.0000272"
Set-TextValue "E12" "  -1.64%  "

# Row 13
Set-TextValue "D13" "9.17"
Set-TextValue "E13" "  -3.32%  "

# Row 14
Set-TextValue "D14" "4.043.70"
Set-TextValue "E14" "  -1.03%  "

# Row 15
Set-TextValue "D15" "3.488.95"
Set-TextValue "E15" "  -0.96%  "

# Row 16
Set-TextValue "E16" "  +0.11%  "

# Row 17
Set-TextValue "D17" "18.38"
Set-TextValue "E17" "  -0.04%  "

# Row 18
Set-TextValue "D18" "12.19"
Set-TextValue "E18" "  +2.34%  "

# Row 19
Set-TextValue "D19" "65.425.02"
Set-TextValue "E19" "  -4.24%  "

# Row 20
Set-TextValue "E20" "  -1.47%  "

# Row 21
Set-TextValue "D21" "413.69"
Set-TextValue "E21" "  +1.13%  "

# Row 22
Set-TextValue "D22" "4.06"
Set-TextValue "E22" "  +2.30%  "

# Row 23
Set-TextValue "D23" "85.67"
Set-TextValue "E23" "  +1.05%  "

# Row 24
Set-TextValue "D24" "4.10"
Set-TextValue "E24" "  -2.76%  "

# Row 25
Set-TextValue "D25" "12.80"
Set-TextValue "E25" "  +7.43%  "

# Row 26
Set-TextValue "D26" "10.79"
Set-TextValue "E26" "  -7.71%  "

# Row 27
Set-TextValue "E27" "  -1.98%  "

# Row 28
Set-TextValue "D28" "9.02"
Set-TextValue "E28" "  +4.50%  "

# Row 29
Set-TextValue "D29" "30.34"
Set-TextValue "E29" "  -0.73%  "

# Row 30
Set-TextValue "D30" "617.27"
Set-TextValue "E30" "  -9.59%  "

# Row 31
Set-TextValue "D31" "6.46"
Set-TextValue "E31" "  -6.22%  "

# Row 32
Set-TextValue "D32" "11.64"
Set-TextValue "E32" "  -0.85%  "

# Row 33
Set-TextValue "E33" "  -1.17%  "

# Row 34
Set-TextValue "D34" "59.42"
Set-TextValue "E34" "  -1.93%  "

# Row 35
Set-TextValue "E35" "  +10.76%  "

# Row 36
Set-TextValue "E36" "  +0.29%  "

# Row 37
Set-TextValue "D37" "0.0₃0790"
Set-TextValue "E37" "  -6.53%  "

# Row 38
Set-TextValue "D38" "37.09"
Set-TextValue "E38" "  -5.18%  "

# Row 39
Set-TextValue "D39" "3.368.16"
Set-TextValue "E39" "  +10.42%  "

# Row 40
Set-TextValue "E40" "  -5.91%  "

# Row 41
Set-TextValue "E41" "  -0.18%  "

# Row 42
Set-TextValue "D42" "3.26"
Set-TextValue "E42" "  -4.22%  "

# Row 43
Set-TextValue "D43" "2.85"
Set-TextValue "E43" "  -5.41%  "

# Row 44
Set-TextValue "D44" "2.78"
Set-TextValue "E44" "  +0.97%  "

# Row 45
Set-TextValue "E45" "  -8.86%  "

# Row 46
Set-TextValue "D46" "0.0415"
Set-TextValue "E46" "  -2.19%  "

# Row 47
Set-TextValue "D47" "3.24"
Set-TextValue "E47" "  +0.84%  "

# Row 48
Set-TextValue "E48" "  +1.69%  "

# Row 49
Set-TextValue "B49" "Monero"
Set-TextValue "C49" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D49" "137.68"
Set-TextValue "E49" "  -0.72%  "

# Row 50
Set-TextValue "B50" "THORChain"
Set-TextValue "C50" "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue "D50" "8.42"
Set-TextValue "E50" "  -5.82%  "

# Row 51
Set-TextValue "D51" "2.87"
Set-TextValue "E51" "  +10.06%  "
